$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '29.421.57'
$ws.Range('D3').Value = '1.850.36'
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('D6').Value = '0.6303'
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '0.07706'
$ws.Range('E8').Value = '  +2.17%  '
$ws.Range('D9').Value = '0.2945'
$ws.Range('E9').Value = '  -0.35%  '
$ws.Range('D10').Value = '24.53'
$ws.Range('E10').Value = '  +0.47%  '
$ws.Range('E11').Value = '  +0.61%  '
$ws.Range('D12').Value = '1.850.74'
$ws.Range('E12').Value = '  -0.89%  '
$ws.Range('D13').Value = '5.030'
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('D14').Value = '0.00001084'
$ws.Range('E14').Value = '  +7.86%  '
$ws.Range('D15').Value = '0.6804'
$ws.Range('D16').Value = '83.76'
$ws.Range('E16').Value = '  +1.09%  '
$ws.Range('D17').Value = '2.108.33'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('D19').Value = '29.440.84'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').Value = '229.42'
$ws.Range('E20').Value = '  +0.97%  '
$ws.Range('E21').Value = '  +0.28%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').Value = '7.457'
$ws.Range('E23').Value = '  -1.09%  '
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('D25').Value = '157.35'
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('D26').Value = '0.1388'
$ws.Range('E26').Value = '  -0.34%  '
$ws.Range('D27').Value = '8.380'
$ws.Range('E27').Value = '  +0.28%  '
$ws.Range('E28').Value = '  +0.31%  '
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('E30').Value = '  +4.87%  '
$ws.Range('D31').Value = '0.05741'
$ws.Range('E31').Value = '  +1.15%  '
$ws.Range('D32').Value = '4.116'
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('D33').Value = '4.051'
$ws.Range('E33').Value = '  +0.86%  '
$ws.Range('D34').Value = '1.851'
$ws.Range('E34').Value = '  +0.41%  '
$ws.Range('E35').Value = '  +0.52%  '
$ws.Range('D36').Value = '0.7099'
$ws.Range('E36').Value = '  -0.31%  '
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('D38').Value = '2.779'
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('D39').Value = '1.229.01'
$ws.Range('E39').Value = '  -2.74%  '
$ws.Range('D40').Value = '0.01799'
$ws.Range('E40').Value = '  -0.89%  '
$ws.Range('D41').Value = '6.459'
$ws.Range('E41').Value = '  +4.11%  '
$ws.Range('D42').Value = '0.9115'
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('D43').Value = '1.000'
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D44').Value = '2.016.94'
$ws.Range('E44').Value = '  -0.74%  '
$ws.Range('D45').Value = '101.88'
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('D46').Value = '66.32'
$ws.Range('E46').Value = '  +0.19%  '
$ws.Range('D47').Value = '0.00000000121'
$ws.Range('E47').Value = '  +3.22%  '
$ws.Range('D48').Value = '7.163'
$ws.Range('E48').Value = '  +1.24%  '
$ws.Range('D49').Value = '0.4024'
$ws.Range('E49').Value = '  -0.44%  '
$ws.Range('D50').Value = '9.021'
$ws.Range('E50').Value = '  -0.62%  '
$ws.Range('E51').Value = '  +0.64%  '
